# Load properties from json
# Renames the {{d2si.*}} placeholders to {{company.*}}, and fills in the
# per-day worked_days / extra_worked_days placeholder cells that were
# previously blank, plus a couple of cosmetic sheet-view tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the d2si.* template placeholders to company.*
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "{{company.logo}}"
$ws.Range("E8").Value = "{{company.employee.name}}"
$ws.Range("E9").Value = "{{company.administrative.name}}"
$ws.Range("E10").Value = "{{company.employee.contract}}"
$ws.Range("H38").Value = "{{company.employee.signature}}"

# ---------------------------------------------------------------------
# 2. Fill in the per-day placeholder cells for the AM/PM worked-days grid
#    (row 17 = AM, row 18 = PM) and its "sum" column (AU).
#    Columns C..AH hold day numbers 1..32 (matching the {{month.day.N}}
#    header in row 16 / row 26).
# ---------------------------------------------------------------------
$firstCol = 3   # column C
$lastCol = 34   # column AH
$sumCol = 47    # column AU

for ($col = $firstCol; $col -le $lastCol; $col++) {
    $day = $col - $firstCol + 1

    $ws.Cells.Item(17, $col).Value = "{{worked_days." + $day + ".AM}}"
    $ws.Cells.Item(18, $col).Value = "{{worked_days." + $day + ".PM}}"

    $ws.Cells.Item(27, $col).Value = "{{extra_worked_days." + $day + ".AM}}"
    $ws.Cells.Item(28, $col).Value = "{{extra_worked_days." + $day + ".PM}}"
    $ws.Cells.Item(29, $col).Value = "{{extra_worked_days." + $day + ".NIGHT}}"
    $ws.Cells.Item(30, $col).Value = "{{extra_worked_days." + $day + ".EARLY_MORNING}}"
}

$ws.Cells.Item(17, $sumCol).Value = "{{worked_days.sum.AM}}"
$ws.Cells.Item(18, $sumCol).Value = "{{worked_days.sum.PM}}"
$ws.Cells.Item(27, $sumCol).Value = "{{extra_worked_days.sum.AM}}"
$ws.Cells.Item(28, $sumCol).Value = "{{extra_worked_days.sum.PM}}"
$ws.Cells.Item(29, $sumCol).Value = "{{extra_worked_days.sum.NIGHT}}"
$ws.Cells.Item(30, $sumCol).Value = "{{extra_worked_days.sum.EARLY_MORNING}}"

# ---------------------------------------------------------------------
# 3. Reset the sheet view / selection back to the top-left corner, and
#    restore the workbook tab-bar ratio.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.TabRatio = 0.677
$ws.Range("A1").Select()
